# Fill in the "Ist Aufwand" (actual effort) row for Phase 1 on the
# "Terminplan" sheet. Row 10 previously only had the "Aktivität/Pflichtenheft
# & Projektplan erstellen" value in column D (1); Ramon now supplies the
# remaining actual-effort figures for the other team members / totals
# columns (E:H), and the row's total (I10, =SUM(D10:H10)) recalculates
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Terminplan")

$ws.Range("E10").Value() = 20
$ws.Range("F10").Value() = 3
$ws.Range("G10").Value() = 0
$ws.Range("H10").Value() = 4
